$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.934.43'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.63%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.320.30'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +7.10%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.66'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '621.00'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.12'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.385'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.317.44'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.785'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.03%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '96.728.51'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.73%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '35.24'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.942.04'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +7.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.49'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.317.22'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +6.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.58'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.03'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '486.09'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +8.35%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.82'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.44%  '
$ws.Range('B23').Value = 'PEPE'
$ws.Range('C23').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000206'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.84%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.65'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '88.01'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.09'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.537.33'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +7.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.240'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.88%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.121'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.26'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.25'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.43'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.152'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '499.77'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.65%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.62'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.66%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.87%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.45%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.29'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.48'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.89%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.785'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +13.76%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '161.02'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.36'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.17%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.51'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.06'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.80%  '
